# Apply the "Office Theme" design colors (the presentation currently uses the
# "Integral" design). The font scheme and format scheme (fills/lines/effects)
# are identical between the two built-in themes already embedded in this
# deck, so re-pointing the 12 theme colors reproduces the new design.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

function RGBFromHex([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$hexColors = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Colors($i).RGB = RGBFromHex $hexColors[$i - 1]
}

# Apply the built-in "Medium Style 2 - Accent 1" table style to the table on
# slide 5 (it previously used the deck's default table style explicitly).
$tableSlide = $p.Slides.Item(5)
for ($i = 1; $i -le $tableSlide.Shapes.Count; $i++) {
    $shape = $tableSlide.Shapes.Item($i)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle("{CFB98ABD-5FD6-494B-832F-7177B4C315BF}")
    }
}
